$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.261.93'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '3.487.77'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''586.77'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').Value = '''134.24'
$ws.Range('E6').Value = '  +1.38%  '
$ws.Range('D7').Value = '3.487.20'
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').Value = '''7.19'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').Value = '''0.376'
$ws.Range('D13').Value = '4.082.26'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '3.496.69'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '64.315.35'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('E18').Value = '  -9.55%  '
$ws.Range('D19').Value = '''9.99'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Value = '''5.70'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').Value = '''13.64'
$ws.Range('E21').Value = '  -5.14%  '
$ws.Range('D22').Value = '''384.34'
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('D24').Value = '3.628.39'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').Value = '''74.15'
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').Value = '''0.0000113'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').Value = '''1.55'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '''0.997'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('D31').Value = '''7.40'
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('D33').Value = '''8.22'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').Value = '3.510.59'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = '''0.148'
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('D37').Value = '''23.36'
$ws.Range('E37').Value = '  -2.26%  '
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('D41').Value = '''161.52'
$ws.Range('E41').Value = '  -4.37%  '
$ws.Range('E42').Value = '  -4.22%  '
$ws.Range('D43').Value = '''0.803'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''25.58'
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '''41.72'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('D47').Value = '''4.38'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = '2.465.62'
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('E51').Value = '  -2.38%  '
